$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1638.5454
$ws.Range("J70").Value = 1662.2
$ws.Range("L70").Value = 4986.6
$ws.Range("N70").Value = -5526.6

$ws.Range("H73").Value = 1638.5454
$ws.Range("J73").Value = 1662.2
$ws.Range("L73").Value = 4986.6
$ws.Range("N73").Value = -6858.6

$ws.Range("H88").Value = 1946.7858
$ws.Range("I88").Value = 1341.5
$ws.Range("J88").Value = 2047.6666
$ws.Range("K88").Value = 1341.5
$ws.Range("L88").Value = 2047.6666
$ws.Range("M88").Value = -935.5
$ws.Range("N88").Value = -2859.6666

$ws.Range("H91").Value = 1946.7858
$ws.Range("I91").Value = 1341.5
$ws.Range("J91").Value = 2047.6666
$ws.Range("K91").Value = 1341.5
$ws.Range("L91").Value = 2047.6666
$ws.Range("M91").Value = 62.5
$ws.Range("N91").Value = -4855.6666

$ws.Range("H137").Value = 3001148.5
$ws.Range("I137").Value = 1220476.8
$ws.Range("J137").Value = 11113098
$ws.Range("K137").Value = 3661430.4
$ws.Range("L137").Value = 33339294
$ws.Range("M137").Value = -3658880.4
$ws.Range("N137").Value = -33344394

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5173.2856
$ws.Range("I2").Value = 1658.25
$ws.Range("J2").Value = 9860
$ws.Range("K2").Value = 1658.25
$ws.Range("L2").Value = 9860
$ws.Range("M2").Value = -1545.25
$ws.Range("N2").Value = -10086

$ws.Range("H45").Value = 4440.9165
$ws.Range("I45").Value = 4087.4285
$ws.Range("J45").Value = 4935.8
$ws.Range("K45").Value = 4087.4285
$ws.Range("L45").Value = 4935.8
$ws.Range("M45").Value = -3710.4285
$ws.Range("N45").Value = -5689.8

$ws.Range("H116").Value = 5173.2856
$ws.Range("I116").Value = 1658.25
$ws.Range("J116").Value = 9860
$ws.Range("K116").Value = 1658.25
$ws.Range("L116").Value = 9860
$ws.Range("M116").Value = 635.75
$ws.Range("N116").Value = -14448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5173.2856
$ws.Range("I3").Value = 1658.25
$ws.Range("J3").Value = 9860
$ws.Range("K3").Value = 1658.25
$ws.Range("L3").Value = 9860
$ws.Range("M3").Value = -1544.25
$ws.Range("N3").Value = -10088

$ws.Range("H87").Value = 10000
$ws.Range("J87").Value = 10000
$ws.Range("L87").Value = 10000
$ws.Range("N87").Value = -12496

$ws.Range("H90").Value = 10000
$ws.Range("J90").Value = 10000
$ws.Range("L90").Value = 30000
$ws.Range("N90").Value = -42480

$ws.Range("H107").Value = 3834.577
$ws.Range("I107").Value = 3994.7646
$ws.Range("J107").Value = 3532
$ws.Range("K107").Value = 3994.7646
$ws.Range("L107").Value = 3532
$ws.Range("M107").Value = -2074.7646
$ws.Range("N107").Value = -7372

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1421.7894
$ws.Range("I16").Value = 1229
$ws.Range("J16").Value = 1636
$ws.Range("K16").Value = 1229
$ws.Range("L16").Value = 1636
$ws.Range("M16").Value = -942
$ws.Range("N16").Value = -2210

$ws.Range("H113").Value = 1421.7894
$ws.Range("I113").Value = 1229
$ws.Range("J113").Value = 1636
$ws.Range("K113").Value = 1229
$ws.Range("L113").Value = 1636
$ws.Range("M113").Value = 941
$ws.Range("N113").Value = -5976

$ws.Range("H122").Value = 1754.8182
$ws.Range("I122").Value = 1107.6154
$ws.Range("J122").Value = 2689.6667
$ws.Range("K122").Value = 3322.8462
$ws.Range("L122").Value = 8069.000100000001
$ws.Range("M122").Value = -872.8462
$ws.Range("N122").Value = -12969.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 142889.42
$ws.Range("I2").Value = 200034.2
$ws.Range("J2").Value = 27.5
$ws.Range("K2").Value = 1200205.2
$ws.Range("L2").Value = 165
$ws.Range("M2").Value = -1200092.2
$ws.Range("N2").Value = -391

$ws.Range("H22").Value = 1646.1538
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 1708.3334
$ws.Range("K22").Value = 2700
$ws.Range("L22").Value = 5125.0002
$ws.Range("M22").Value = -2531
$ws.Range("N22").Value = -5463.0002

$ws.Range("H27").Value = 1646.1538
$ws.Range("I27").Value = 900
$ws.Range("J27").Value = 1708.3334
$ws.Range("K27").Value = 2700
$ws.Range("L27").Value = 5125.0002
$ws.Range("M27").Value = -2598
$ws.Range("N27").Value = -5329.0002

$ws.Range("H49").Value = 7833
$ws.Range("J49").Value = 8499.6
$ws.Range("L49").Value = 25498.8
$ws.Range("N49").Value = -25810.8

$ws.Range("H75").Value = 1753.75
$ws.Range("I75").Value = 500
$ws.Range("J75").Value = 2171.6667
$ws.Range("K75").Value = 1500
$ws.Range("L75").Value = 6515.000100000001
$ws.Range("M75").Value = -502
$ws.Range("N75").Value = -8511.000100000001

$ws.Range("H78").Value = 1753.75
$ws.Range("I78").Value = 500
$ws.Range("J78").Value = 2171.6667
$ws.Range("K78").Value = 4500
$ws.Range("L78").Value = 19545.0003
$ws.Range("M78").Value = 492
$ws.Range("N78").Value = -29529.0003

$ws.Range("H87").Value = 12239.25
$ws.Range("I87").Value = 9478.5
$ws.Range("K87").Value = 28435.5
$ws.Range("M87").Value = -27187.5

$ws.Range("H90").Value = 12239.25
$ws.Range("I90").Value = 9478.5
$ws.Range("K90").Value = 85306.5
$ws.Range("M90").Value = -79066.5

$ws.Range("H113").Value = 668.86206
$ws.Range("I113").Value = 602.6429000000001
$ws.Range("J113").Value = 730.6667
$ws.Range("K113").Value = 1807.9287
$ws.Range("L113").Value = 2192.0001
$ws.Range("M113").Value = 362.0712999999998
$ws.Range("N113").Value = -6532.0001

$ws.Range("H122").Value = 2898.55
$ws.Range("I122").Value = 320.42856
$ws.Range("J122").Value = 3318.2441
$ws.Range("K122").Value = 2883.85704
$ws.Range("L122").Value = 29864.1969
$ws.Range("M122").Value = -433.8570399999999
$ws.Range("N122").Value = -34764.1969

$ws.Range("H132").Value = 1197263
$ws.Range("I132").Value = 1316689.2
$ws.Range("K132").Value = 11850202.8
$ws.Range("M132").Value = -11847672.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 38088.5
$ws.Range("J82").Value = 38088.5
$ws.Range("L82").Value = 38088.5
$ws.Range("N82").Value = -38854.5

$ws.Range("H85").Value = 38088.5
$ws.Range("J85").Value = 38088.5
$ws.Range("L85").Value = 38088.5
$ws.Range("N85").Value = -40740.5

$ws.Range("H96").Value = 1666.6666
$ws.Range("I96").Value = 1740
$ws.Range("J96").Value = 1740
$ws.Range("K96").Value = 1740
$ws.Range("L96").Value = 1740
$ws.Range("M96").Value = -367
